$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Remove the two old rows (18 and 27) that documented the "404C2"/"412C"
# errors. Clearing their contents removes them from the saved sheetData
# (and prunes the now-unused shared strings) without shifting other rows.
# ---------------------------------------------------------------------------
$ws.Range("A18:E18").ClearContents() | Out-Null
$ws.Range("A27:E27").ClearContents() | Out-Null

# ---------------------------------------------------------------------------
# Add the new row 12 describing the "401A" fix (2014 is now completely
# clean, this was the last remaining error to resolve).
# ---------------------------------------------------------------------------
$ws.Range("A12").Value = 2014
$ws.Range("B12").Value = "401A"
$ws.Range("C12").Value = 28
$ws.Range("D12").Value = "Month is 24, day is 3 theyre flipped"
$ws.Range("E12").Value = "make month 3 and day 24"

# ---------------------------------------------------------------------------
# Add the new verification row 26 (L26 date, M26 = L26+27 formula), both
# formatted with the built-in short-date number format (numFmtId 14).
# ---------------------------------------------------------------------------
$ws.Range("L26").Value = 41720
$ws.Range("L26").NumberFormat = "mm-dd-yy"
$ws.Range("M26").Formula = "=L26+27"
$ws.Range("M26").NumberFormat = "mm-dd-yy"

# ---------------------------------------------------------------------------
# Update the active cell selection to match the author's saved view.
# ---------------------------------------------------------------------------
$ws.Range("E16").Select() | Out-Null
